# "update New Data Process code"
# Insert a new "Stop Gen" data row (400/400/400) right after the first
# case's "Average Time" row (was row 4, pushing all subsequent rows down
# by one -- old row 42/43 become 43/44, etc.), then leave the selection
# where the author's cursor ended up (D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 4 (shifts rows 4.. down to 5..)
$ws.Rows("4:4").Insert()

# Fill in the new row's values
$ws.Range("A4").Value = "Stop Gen"
$ws.Range("B4").Value = 400
$ws.Range("C4").Value = 400
$ws.Range("D4").Value = 400

# Match the author's final cell selection
$ws.Range("D8").Select() | Out-Null
